$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing D:K data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy cell formatting (number formats/styles) from column F (the old column D, now shifted)
# into the two newly inserted columns D and E, for the full data range.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newValues = @{
    7 = @(43465, 43373)
    8 = @(741600, 655400)
    9 = @(309800, 305000)
    10 = @(431800, 350400)
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(26200, 7300)
    15 = @(0, 0)
    17 = @(709400, 653100)
    18 = @(32200, 2300)
    20 = @(-5800, -1500)
    21 = @(84200, 40900)
    22 = @(51400, 50200)
    23 = @(-25000, -49400)
    24 = @(53000, -38700)
    25 = @(0, 0)
    26 = @(-78000, -10700)
    27 = @(-78000, -10700)
    28 = @(0, 0)
    29 = @(7700, -400)
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(5800, 1500)
    33 = @(-70300, -11100)
    34 = @(0, 0)
    35 = @(-70300, -11100)
    38 = @(43465, 43373)
    41 = @(87300, 61800)
    42 = @(0, 0)
    43 = @(431300, 441600)
    44 = @(523200, 587600)
    45 = @(152000, 171200)
    46 = @(1193800, 1262200)
    47 = @(0, 0)
    48 = @(354500, 358400)
    49 = @(1205900, 1252200)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(262600, 315500)
    53 = @(0, 0)
    54 = @(3016800, 3188300)
    57 = @(332100, 363200)
    58 = @(357400, 441700)
    59 = @(430900, 411600)
    60 = @(1120400, 1216500)
    61 = @(2727700, 2731400)
    62 = @(225500, 228600)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(4073600, 4176500)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(-1855000, -1784700)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(-1056800, -988200)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(-70300, -11100)
    83 = @(57800, 40100)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(125900, -106600)
    91 = @(-15600, -11800)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-15600, -11800)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-83500, 99700)
    101 = @(-1800, -1100)
    102 = @(25000, -19800)
}


foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}

Write-Output "done"
